$wb = $excel.ActiveWorkbook

# --- Sheet: Summary ---
$ws = $wb.Worksheets.Item("Summary")
$ws.Range("B2").Value = 0.0498220640569395
$ws.Range("C2").Value = 0.0498220640569395
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 0.09491525423728814
$ws.Range("F2").Value = 0.2077151335311573
$ws.Range("G2").Value = 0.5768621236133122
$ws.Range("H2").Value = 0.8104935794542536
$ws.Range("I2").Value = 28
$ws.Range("J2").Value = 534
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0

# --- Sheet: Classification Report ---
$ws = $wb.Worksheets.Item("Classification Report")

# Row 2 (class "0")
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0

# Row 3 (class "1")
$ws.Range("B3").Value = 0.0498220640569395
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 0.09491525423728814

# Row 4 (accuracy)
$ws.Range("B4").Value = 0.0498220640569395
$ws.Range("C4").Value = 0.0498220640569395
$ws.Range("D4").Value = 0.0498220640569395
$ws.Range("E4").Value = 0.0498220640569395

# Row 5 (macro avg)
$ws.Range("B5").Value = 0.02491103202846975
$ws.Range("C5").Value = 0.5
$ws.Range("D5").Value = 0.04745762711864407

# Row 6 (weighted avg)
$ws.Range("B6").Value = 0.002482238066893783
$ws.Range("C6").Value = 0.0498220640569395
$ws.Range("D6").Value = 0.004728873876590867

# --- Sheet: Confusion Matrix ---
$ws = $wb.Worksheets.Item("Confusion Matrix")

# Row 2 (Actual 0)
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 534

# Row 3 (Actual 1)
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 28
